# Standardize schema headers to snake_case (sheet-by-sheet rename), and
# restructure SpeciesRegisterEntry / OrganizationRegisterEntry to the new
# snake_case column layout.

$wb = $excel.ActiveWorkbook

# NOTE: this PS host does not resolve named (-Param value) arguments on
# user functions, so Set-HeaderRow is called positionally everywhere below.
function Set-HeaderRow {
    param($Worksheet, $Headers)
    for ($i = 0; $i -lt $Headers.Length; $i++) {
        $Worksheet.Cells.Item(1, $i + 1).Value = $Headers[$i]
    }
}

# Organization: camelCase -> snake_case (same column layout)
$ws = $wb.Worksheets.Item("Organization")
Set-HeaderRow $ws @(
    "name", "domain", "site_type_ids", "species_ids", "url_path", "internal_path",
    "slug", "id", "created_at", "created_by_id", "updated_at", "updated_by_id", "organization_id"
)

# Site: camelCase -> snake_case (same column layout)
$ws = $wb.Worksheets.Item("Site")
Set-HeaderRow $ws @(
    "site_type_id", "name", "group_id_hierarchy", "description", "url_path", "internal_path",
    "slug", "id", "created_at", "created_by_id", "updated_at", "updated_by_id", "organization_id"
)

# Location: camelCase -> snake_case (same column layout)
$ws = $wb.Worksheets.Item("Location")
Set-HeaderRow $ws @(
    "group_type_id", "name", "site_id", "parent_id", "description", "capacity",
    "url_path", "internal_path", "slug", "id", "created_at", "created_by_id",
    "updated_at", "updated_by_id", "organization_id"
)

# Genet: camelCase -> snake_case (same column layout)
$ws = $wb.Worksheets.Item("Genet")
Set-HeaderRow $ws @(
    "name", "species_id", "genet_type_id", "sf_id", "clonal_id", "accession_number",
    "url_path", "internal_path", "slug", "id", "created_at", "created_by_id",
    "updated_at", "updated_by_id", "organization_id"
)

# Coral: camelCase -> snake_case (same column layout)
$ws = $wb.Worksheets.Item("Coral")
Set-HeaderRow $ws @(
    "name", "genet_id", "species_id", "site_id", "group_id", "coral_type_id",
    "quantity", "coral_size", "url_path", "internal_path", "slug", "id",
    "created_at", "created_by_id", "updated_at", "updated_by_id", "organization_id"
)

# Species: unchanged

# Person: camelCase -> snake_case (same column layout)
$ws = $wb.Worksheets.Item("Person")
Set-HeaderRow $ws @(
    "name", "email", "image_url", "id", "created_at", "created_by_id",
    "updated_at", "updated_by_id", "organization_id"
)

# Event: camelCase -> snake_case (same column layout)
$ws = $wb.Worksheets.Item("Event")
Set-HeaderRow $ws @(
    "event_type_id", "record_id", "record_model_type", "url_path", "internal_path",
    "slug", "id", "created_at", "created_by_id", "updated_at", "updated_by_id", "organization_id"
)

# SpeciesRegisterEntry: restructured layout (A1:H1 -> A1:I1)
$ws = $wb.Worksheets.Item("SpeciesRegisterEntry")
Set-HeaderRow $ws @(
    "id", "created_at", "updated_at", "common_name", "genus", "specific_epithet",
    "scientific_name", "photo_url", "tags"
)

# OrganizationRegisterEntry: restructured layout (A1:E1 -> A1:J1)
$ws = $wb.Worksheets.Item("OrganizationRegisterEntry")
Set-HeaderRow $ws @(
    "id", "created_at", "updated_at", "name", "description", "region",
    "website_url", "contact_email", "logo_url", "is_active"
)
